# Update DG pictures to be consistent in Add Edit Delete List
#
# This script:
#  1. Updates the fixed "date last saved" footer field (10/22/19 -> 10/24/2019)
#     on the slide master and every slide layout.
#  2. Repositions the big background rounded rectangle slightly.
#  3. Swaps the DeleteBondCommand / EditBondCommand rounded-rectangle boxes so
#     DeleteBondCommand sits where EditBondCommand used to be (and vice versa),
#     nudging their y-offsets very slightly.
#  4. Updates the three connector lines that drop down from the
#     {abstract}Command bar into the Add/Delete/Edit boxes so their geometry
#     matches the new box layout.

# ---------------------------------------------------------------------------
# Helper: PowerPoint's Shape.Left/Top/Width/Height are single-precision
# (float32) in the real object model, and the host truncates when it turns
# the point value back into EMU (floor(float32(pt) * 12700)). Nudge the
# point value up by the smallest amount needed so the stored EMU lands on
# the exact target instead of one EMU short.
# ---------------------------------------------------------------------------
function Emu-ToSafePt($targetEmu) {
    $emuPerPt = 12700.0
    $base = $targetEmu / $emuPerPt
    for ($i = 0; $i -le 4000; $i++) {
        $cand = $base + ($i * 0.0000001)
        $f32 = [float]$cand
        $emu = [math]::Floor([double]$f32 * $emuPerPt)
        if ($emu -eq $targetEmu) {
            return $cand
        }
    }
    return $base
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Fixed date field: 10/22/19 -> 10/24/2019 (slide master + all layouts)
# ---------------------------------------------------------------------------

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePh = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDatePh = $true
            }
        } catch {
        }
        if ($isDatePh -and $sh.HasTextFrame -eq -1) {
            if ($sh.TextFrame.TextRange.Text -eq "10/22/19") {
                $sh.TextFrame.TextRange.Text = "10/24/2019"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $cl = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $cl.Shapes
}

# ---------------------------------------------------------------------------
# Locate the slide + the "Group 3" group that holds the diagram shapes.
# ---------------------------------------------------------------------------

$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(1)
$items = $grp.GroupItems

function Get-ShapeById($items, $id) {
    for ($i = 1; $i -le $items.Count; $i++) {
        $sh = $items.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 2. Background rounded rectangle nudges down slightly.
# ---------------------------------------------------------------------------

$bigRect = Get-ShapeById $items 5
$bigRect.Top = Emu-ToSafePt 292903

# ---------------------------------------------------------------------------
# 3. Swap DeleteBondCommand (id 11) and EditBondCommand (id 12) boxes.
# ---------------------------------------------------------------------------

$deleteBox = Get-ShapeById $items 11
$editBox = Get-ShapeById $items 12

$deleteBox.Left = Emu-ToSafePt 6785800
$deleteBox.Top = Emu-ToSafePt 2610170

$editBox.Left = Emu-ToSafePt 4430940
$editBox.Top = Emu-ToSafePt 2610171

# ---------------------------------------------------------------------------
# 4. Update the connector lines feeding into the Add/Delete/Edit boxes.
# ---------------------------------------------------------------------------

# Straight Connector 16 (id 17) -> starts at AddBondCommand box (id 10)
$conn16 = Get-ShapeById $items 17
$conn16.Left = Emu-ToSafePt 3034684
$conn16.Top = Emu-ToSafePt 2247003
$conn16.Width = Emu-ToSafePt 1569
$conn16.Height = Emu-ToSafePt 413048

# Straight Connector 17 (id 18) -> starts at DeleteBondCommand box (id 11)
$conn17 = Get-ShapeById $items 18
$conn17.Left = Emu-ToSafePt 7777144
$conn17.Top = Emu-ToSafePt 2247003
$conn17.Width = Emu-ToSafePt 0
$conn17.Height = Emu-ToSafePt 363167

# Straight Connector 18 (id 19) -> starts at EditBondCommand box (id 12)
$conn18 = Get-ShapeById $items 19
$conn18.Left = Emu-ToSafePt 5422284
$conn18.Top = Emu-ToSafePt 2247003
$conn18.Width = Emu-ToSafePt 0
$conn18.Height = Emu-ToSafePt 363168
